# Updated symbol list on Thu Dec 15 18:27:49 UTC 2022 with GitHub Actions
#
# Applies the per-row "Price" (column D) refreshes plus the BKEXToken /
# CEJI / KickToken re-ranking (rows 41-43) exactly as captured by the
# source diff. Column D holds numeric-looking values that are stored as
# *text* in the workbook, so every D-cell write goes through
# Set-TextValue, which forces the General/@ -> Normal dance needed to
# keep Excel from silently re-typing the cell as a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- simple Price (column D) updates -------------------------------------
Set-TextValue "D2"  "261.19"
Set-TextValue "D4"  "6.211"
Set-TextValue "D5"  "0.06092"
Set-TextValue "D7"  "6.706"
Set-TextValue "D8"  "1.359"
Set-TextValue "D9"  "0.7993"
Set-TextValue "D10" "0.1574"
Set-TextValue "D11" "0.08131"
Set-TextValue "D12" "0.03324"
Set-TextValue "D13" "0.03137"
Set-TextValue "D14" "0.09268"
Set-TextValue "D15" "3.893"
Set-TextValue "D16" "0.001695"
Set-TextValue "D17" "0.04825"
Set-TextValue "D18" "0.0006215"
Set-TextValue "D19" "0.006224"
Set-TextValue "D20" "0.001103"
Set-TextValue "D21" "0.003374"
Set-TextValue "D25" "0.3372"
Set-TextValue "D40" "0.04609"

# --- rows 41-43 re-ranked: BKEXToken / CEJI / KickToken -------------------
# New order: KickToken, BKEXToken, CEJI (with refreshed prices + the
# "<rank>Name Ticker" Volume(1h) strings bumped to match their new rows).
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.007205"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1119"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003132"
$ws.Range("E43").Value = "42CEJICEJI"

# --- remaining simple Price (column D) updates ----------------------------
Set-TextValue "D44" "0.01021"
Set-TextValue "D46" "0.00006027"
Set-TextValue "D49" "0.05404"
